$d = $word.ActiveDocument
$ellipsis = [char]0x2026

# 1. Update the "I believe understanding..." blurb about math/functional programming.
$d.Content.Find.Execute(
    "I believe understanding math and functional programming is essential for one to excel at writing sustainable software.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I believe understanding mathematics is essential for one to excel at writing sustainable software and mathematics is made manifest in functional programming concepts.",
    2)

# 2. Replace the terse "flow, jest, ESLint" bullet with the "code quality" bullet.
$d.Content.Find.Execute(
    "flow, jest, ESLint",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "code quality (ESLint, flow, etc" + $ellipsis + ")",
    2)

# 3. Replace the "I pretty much am familiar..." bullet with the testing-tools bullet.
$d.Content.Find.Execute(
    "I pretty much am familiar with everything" + $ellipsis + "not that I am a master of everything, but that I have at least tried everything and understand the purpose all items in the ecosystem.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "unit/integration/fuzz/acceptance testing (jest, AVA, mocha, jasmine, etc" + $ellipsis + ")",
    2)

# 4. Tweak the "I keep track of current trends..." sentence.
$d.Content.Find.Execute(
    "I keep track of current trends and developments" + $ellipsis + "it is kind of a hobby",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I follow track of current trends and developments" + $ellipsis + "it is one of my hobbies :)",
    2)

# 5. Shorten the bold "example" lead-in to "ex" (match whole word so the
#    "examples" inside the earlier "for examples)" text is left untouched).
$d.Content.Find.Execute(
    "example",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "ex",
    2)
